# Auto-generated script to apply betting-odds value updates for 2025-12-31 workbook
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F2").Value = 3.95
$ws.Range("G2").Value = 4
$ws.Range("I2").Value = 2.12
$ws.Range("J2").Value = 3.65
$ws.Range("K2").Value = 3.7
$ws.Range("L2").Value = 1.39
$ws.Range("N2").Value = 4
$ws.Range("O2").Value = 1.31
$ws.Range("P2").Value = 2.02
$ws.Range("Q2").Value = 1.97
$ws.Range("S2").Value = 3.4
$ws.Range("T2").Value = 1.78
$ws.Range("U2").Value = 2.2
$ws.Range("V2").Value = 1.89
$ws.Range("W2").Value = 1.33
$ws.Range("AB2").Value = 15
$ws.Range("AD2").Value = 10.5
$ws.Range("AF2").Value = 27
$ws.Range("AH2").Value = 19.5
$ws.Range("AI2").Value = 34
$ws.Range("AK2").Value = 46
$ws.Range("AM2").Value = 85
$ws.Range("AN2").Value = 44
$ws.Range("AO2").Value = 15
$ws.Range("G3").Value = 1.3
$ws.Range("K3").Value = 7.4
$ws.Range("N3").Value = 5
$ws.Range("O3").Value = 1.22
$ws.Range("P3").Value = 2.22
$ws.Range("T3").Value = 1.04
$ws.Range("U3").Value = 1.04
$ws.Range("W3").Value = 4.3
$ws.Range("AB3").Value = 11
$ws.Range("AC3").Value = 17
$ws.Range("AD3").Value = 150
$ws.Range("AN3").Value = 5
$ws.Range("F4").Value = 4.5
$ws.Range("H4").Value = 1.68
$ws.Range("I4").Value = 1.82
$ws.Range("K4").Value = 4.9
$ws.Range("L4").Value = 1.32
$ws.Range("N4").Value = 4.2
$ws.Range("P4").Value = 2.22
$ws.Range("Q4").Value = 1.61
$ws.Range("R4").Value = 1.45
$ws.Range("S4").Value = 2.46
$ws.Range("T4").Value = 1.65
$ws.Range("U4").Value = 2.12
$ws.Range("V4").Value = 2.18
$ws.Range("W4").Value = 1.23
$ws.Range("Y4").Value = 990
$ws.Range("Z4").Value = 23
$ws.Range("AB4").Value = 990
$ws.Range("AD4").Value = 18
$ws.Range("AH4").Value = 40
$ws.Range("AO4").Value = 38
$ws.Range("F5").Value = 3.9
$ws.Range("H5").Value = 1.91
$ws.Range("L5").Value = 1.38
$ws.Range("P5").Value = 1.89
$ws.Range("R5").Value = 1.35
$ws.Range("W5").Value = 1.28
$ws.Range("Z5").Value = 13.5
$ws.Range("AF5").Value = 34
$ws.Range("AJ5").Value = 1000
$ws.Range("G6").Value = 4.5
$ws.Range("H6").Value = 1.89
$ws.Range("I6").Value = 2.22
$ws.Range("J6").Value = 3.65
$ws.Range("K6").Value = 4.2
$ws.Range("L6").Value = 1.3
$ws.Range("N6").Value = 3.65
$ws.Range("Q6").Value = 1.58
$ws.Range("T6").Value = 1.7
$ws.Range("U6").Value = 2.4
$ws.Range("AC6").Value = 42
$ws.Range("G7").Value = 9.4
$ws.Range("I7").Value = 1.5
$ws.Range("J7").Value = 4.5
$ws.Range("K7").Value = 5.4
$ws.Range("N7").Value = 4.6
$ws.Range("U7").Value = 1.92
$ws.Range("V7").Value = 2.96
$ws.Range("Y7").Value = 980
$ws.Range("Z7").Value = 1000
$ws.Range("AC7").Value = 42
$ws.Range("F8").Value = 3.35
$ws.Range("G8").Value = 3.6
$ws.Range("I8").Value = 2.1
$ws.Range("J8").Value = 4.1
$ws.Range("V8").Value = 1.9
$ws.Range("W8").Value = 1.38
$ws.Range("Y8").Value = 16
$ws.Range("AB8").Value = 23
$ws.Range("AI8").Value = 26
$ws.Range("AK8").Value = 34
$ws.Range("H10").Value = 1.65
$ws.Range("I10").Value = 1.67
$ws.Range("J10").Value = 4.2
$ws.Range("L10").Value = 1.3
$ws.Range("N10").Value = 4
$ws.Range("P10").Value = 2
$ws.Range("Q10").Value = 1.82
$ws.Range("S10").Value = 3.1
$ws.Range("U10").Value = 1.99
$ws.Range("V10").Value = 2.48
$ws.Range("X10").Value = 16.5
$ws.Range("Y10").Value = 9.4
$ws.Range("AA10").Value = 17
$ws.Range("AE10").Value = 18.5
$ws.Range("AF10").Value = 48
$ws.Range("AN10").Value = 90
$ws.Range("AO10").Value = 9.800000000000001
$ws.Range("F11").Value = 13.5
$ws.Range("G11").Value = 17.5
$ws.Range("I11").Value = 1.26
$ws.Range("L11").Value = 1.17
$ws.Range("N11").Value = 8.800000000000001
$ws.Range("O11").Value = 1.11
$ws.Range("R11").Value = 2.04
$ws.Range("S11").Value = 1.79
$ws.Range("V11").Value = 4.8
$ws.Range("W11").Value = 1.06
$ws.Range("Y11").Value = 16.5
$ws.Range("Z11").Value = 11
$ws.Range("AA11").Value = 11
$ws.Range("AB11").Value = 80
$ws.Range("AC11").Value = 19.5
$ws.Range("AH11").Value = 32
$ws.Range("AI11").Value = 32
$ws.Range("AN11").Value = 160
$ws.Range("F12").Value = 4.9
$ws.Range("H12").Value = 1.64
$ws.Range("I12").Value = 1.73
$ws.Range("J12").Value = 4.2
$ws.Range("K12").Value = 5.4
$ws.Range("L12").Value = 1.27
$ws.Range("N12").Value = 4.9
$ws.Range("O12").Value = 1.2
$ws.Range("P12").Value = 2.34
$ws.Range("R12").Value = 1.53
$ws.Range("S12").Value = 2.36
$ws.Range("T12").Value = 1.62
$ws.Range("U12").Value = 2.16
$ws.Range("V12").Value = 2.36
$ws.Range("W12").Value = 1.22
$ws.Range("Y12").Value = 12.5
$ws.Range("AC12").Value = 12
$ws.Range("AG12").Value = 22
$ws.Range("AK12").Value = 150
$ws.Range("AM12").Value = 580
$ws.Range("AO12").Value = 8
$ws.Range("G13").Value = 1.32
$ws.Range("H13").Value = 11.5
$ws.Range("J13").Value = 6
$ws.Range("L13").Value = 1.18
$ws.Range("M13").Value = 1.03
$ws.Range("Q13").Value = 1.55
$ws.Range("R13").Value = 1.66
$ws.Range("S13").Value = 2.22
$ws.Range("T13").Value = 1.9
